# Add "Psychologist Appointment" and "Adherence Counselor Appointment" to the
# "appoint" choice list in the "choices" sheet of the appointment form.
#
# The existing "internal referral" / "external referral" choices (rows 6-7)
# are pushed down to make room for the two new choices, and the "lab" choices
# (cd4 count / viral load, previously rows 8-9) shift down accordingly to
# rows 10-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Insert two blank rows before the current row 8 ("lab" / "cd4 count" ...).
# This pushes the existing rows 8-9 down to 10-11, preserving their
# formatting (style, number formats, etc.).
$ws.Rows("8:9").Insert()

# Turn the (now duplicated) rows 6-7 into the two new appointment types.
$ws.Range("A6").Value = "appoint"
$ws.Range("B6").Value = "psychologist appointment"
$ws.Range("C6").Value = "Psychologist Appointment"

$ws.Range("A7").Value = "appoint"
$ws.Range("B7").Value = "adherence counselor appointment"
$ws.Range("C7").Value = "Adherence Counselor Appointment"

# Re-create the "internal referral" / "external referral" choices that used
# to live in rows 6-7 in the newly inserted rows 8-9.
$ws.Range("A8").Value = "appoint"
$ws.Range("B8").Value = "internal referral"
$ws.Range("C8").Value = "Internal Referral "

$ws.Range("A9").Value = "appoint"
$ws.Range("B9").Value = "external referral"
$ws.Range("C9").Value = "External Referral "

# Widen columns B and C so the new, longer choice labels are fully visible.
# (Target widths are 36.25 / 38.63 "characters"; Excel quantizes ColumnWidth
# to a whole number of pixels internally, so these inputs are chosen to land
# on the closest representable width to the target.)
$ws.Columns("B").ColumnWidth = 35.25
$ws.Columns("C").ColumnWidth = 37.75
